# finished round 3 of cellconc in vert vs horz
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new (mostly blank) row above row 36, shifting the existing ---
# --- rows 36.. down by one.                                             ---
$ws.Rows("36:36").Insert()
$ws.Rows("36:36").Clear()

# --- Append four new data rows (57-60) at the bottom of the table.       ---
# --- Numeric columns first (order doesn't affect shared-string table).   ---
$ws.Range("B57").Value = 5
$ws.Range("C57").Value = "V"
$ws.Range("D57").Value = 610.6
$ws.Range("E57").Value = 4.1073000000000004
$ws.Range("F57").Value = 1198
$ws.Range("G57").Value = 196.7
$ws.Range("H57").Value = 2508

$ws.Range("B58").Value = 5
$ws.Range("C58").Value = "V"
$ws.Range("D58").Value = 596.29999999999995
$ws.Range("E58").Value = 4.1288999999999998
$ws.Range("F58").Value = 1198
$ws.Range("G58").Value = 191.1
$ws.Range("H58").Value = 2462

$ws.Range("B59").Value = 5
$ws.Range("C59").Value = "V"
$ws.Range("D59").Value = 598.79999999999995
$ws.Range("E59").Value = 4.1215000000000002
$ws.Range("F59").Value = 1198
$ws.Range("G59").Value = 192.3
$ws.Range("H59").Value = 2468

$ws.Range("B60").Value = 5
$ws.Range("C60").Value = "V"
$ws.Range("D60").Value = 576.70000000000005
$ws.Range("E60").Value = 4.1475
$ws.Range("F60").Value = 1198
$ws.Range("G60").Value = 186.5
$ws.Range("H60").Value = 2392

# --- Text columns, written in the same order the lab notebook entries ---
# --- were originally authored so the shared-string table lines up.   ---
$ws.Range("A57").Value = "D20151105T221748"
$ws.Range("I57").Value = "lots of missed rois, ypos when way high halfway out of FOV"
$ws.Range("J57").Value = "no clumps of junk, cellconc use all signals, does it calculate conc including zerosize rois?"

$ws.Range("A58").Value = "D20151105T224008"
$ws.Range("I58").Value = "still lots of missed rois, high ypos off FOV ~450-1030"
$ws.Range("J58").Value = "no junk or visible bubbles"
$ws.Range("K58").Value = "from pmtAvsB most of missed rois show signals that they're beads and just out of camera sight"

$ws.Range("I59").Value = "ypos still high off top but not as many missed, 300-1030"
$ws.Range("J59").Value = "no junk, slightly more bead doublets"
$ws.Range("A59").Value = "D20151105T230227"
$ws.Range("K59").Value = "length(tind)=26 which is really weird because a lot of very high ypos"

$ws.Range("L58").Value = "length(tind)=310"
$ws.Range("K57").Value = "length(tind)=782"

$ws.Range("A60").Value = "D20151105T232446"
$ws.Range("I60").Value = "ypos high but core completely in FOV, a few tiny pieces of junk and doublets and triplets but no large clumps, seems acceptable"
$ws.Range("J60").Value = "length(tind)=12"

$ws.Range("I36").Value = "fresh bead in FSW just made to start over"

# --- Restore selection to match the saved view state ---
$ws.Range("I37").Select()
